$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("B1").Value = "plan_price_monthly"
$ws.Range("C1").Value = "plan_price_semi_annual"
$ws.Range("D1").Value = "plan_price_annual"
$ws.Range("E1").Value = "plan_duration"
$ws.Range("F1").Value = "subscribable"
$ws.Range("G1").Value = "visible"
$ws.Range("H1").Value = "opspi_account_id"

# --- Row 2 : hsphere_plan_10 ---
$ws.Range("A2").Value = "hsphere_plan_10"
$ws.Range("B2").Value = 100.0
$ws.Range("C2").Value = 500.0
$ws.Range("D2").Value = 1000.0
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = 1.0

# --- Row 3 : hsphere_plan_11 ---
$ws.Range("A3").Value = "hsphere_plan_11"
$ws.Range("B3").Value = 50.0
$ws.Range("C3").Value = 300.0
$ws.Range("D3").Value = 600.0
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = $true
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = 1.0

# --- Row 4 : hsphere_plan_12 (new) ---
$ws.Range("A4").Value = "hsphere_plan_12"
$ws.Range("B4").Value = 20.0
$ws.Range("C4").Value = 150.0
$ws.Range("D4").Value = 300.0
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = $true
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 1.0

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 17.88
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 18.75
$ws.Columns.Item(4).ColumnWidth = 17.5
$ws.Columns.Item(6).ColumnWidth = $null
$ws.Columns.Item(8).ColumnWidth = 15.13
